# Start Year Number of EV Chargers.xlsx
#
# "updates to eu files to match new us structure"
#
# The SYNoEVC sheet previously carried a full 2021-2050 projection (pulled in
# part from an external workbook via `=[1]Calculations!B2`). The new
# structure trims this down to just the 2021 start-year figure (a literal
# value, no external link) plus a single trailing blank, styled cell - i.e.
# the same shape used by the US input-data files.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYNoEVC")

# Row 1 (year headers): keep A1 ("Number of chargers") and B1 (2021), drop
# the rest of the year columns (C1:AE1 -> 2022..2050).
$ws.Range("C1:AE1").Clear()

# Row 2 (data): B2 becomes a plain literal (the external-workbook formula
# `=[1]Calculations!B2` is replaced by its resolved value) and every column
# after it is dropped, except C2 which survives as an empty, styled cell
# (matching the trailing blank cell pattern from the US files).
$ws.Range("B2").Value = 299178
$ws.Range("D2:AF2").Clear()
$ws.Range("C2").ClearContents()

# Selection moves from C2 to B8 to match the refreshed sheet view.
$ws.Range("B8").Select()
